# MasterTable_process_a.xlsx - corrections and configuration tables for optimizer
#
# Changes applied:
#  1. Column F ("USE_ACTUAL_MODEL") values for rows 2-8 are renamed from
#     "MLA" to "PR_A_Y1" (the shared string "MLA" is replaced throughout
#     the table by the new label "PR_A_Y1").
#  2. The active selection on the sheet moves from B14 to G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USE_ACTUAL_MODEL column (F2:F8) to the new value "PR_A_Y1"
$ws.Range("F2:F8").Value = "PR_A_Y1"

# Update the selected/active cell shown when the sheet is opened
$ws.Range("G7").Select()
